$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial for every data row (row 2 onward).
# Determine the last used row in column C dynamically, then bump every value
# in that column by one day (e.g. 46060 -> 46061).
$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End($xlUp).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}
